$d = $word.ActiveDocument

# Update the version number shown in the document body ("2021.1-IT" -> "2020.2-IT")
$d.Content.Find.Execute("2021.1", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2020.2", 2)

# Update the version number shown in the footer ("Versione 2021.1-IT" -> "Versione 2020.2-IT")
foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute("2021.1", $true, $false, $false, $false, $false, `
                                     $true, 1, $false, "2020.2", 2)
        }
    }
}
